$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation (2021-11-15, Volumen=160) is inserted as the new
# row 10; every existing row from 10 down shifts one row lower (old row 208
# becomes row 209). Insert a whole row so formatting/styles shift correctly.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the same fixed attributes as every
# other Ciboulette/Femacal de La Calera record, but with the new date and
# volume values.
$ws.Cells.Item(10, 1).Value = 3
$ws.Cells.Item(10, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(10, 3).Value = "Coquimbo"
$ws.Cells.Item(10, 4).Value = 44515
$ws.Cells.Item(10, 5).Value = 5
$ws.Cells.Item(10, 6).Value = 100112039
$ws.Cells.Item(10, 7).Value = "Ciboulette"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 160
$ws.Cells.Item(10, 11).Value = 1500
$ws.Cells.Item(10, 12).Value = 1500
$ws.Cells.Item(10, 13).Value = 1500
$ws.Cells.Item(10, 14).Value = "`$/docena de atados"
$ws.Cells.Item(10, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(10, 16).Value = 500
$ws.Cells.Item(10, 17).Value = 3
$ws.Cells.Item(10, 18).Value = "Hortaliza"
